# "fixed inputs for power plants"
# - Sheet1: rows that referenced the "traderes" source now reference "emlab"
# - Sheet2 ("from traderes to emlab") renamed to "data per year" and filled
#   with the per-year market/power-plant data pulled from competes-emlab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: fix the data source for a few inputs -------------------------
$ws1.Range("B2").Value = "emlab"
$ws1.Range("B3").Value = "emlab"
$ws1.Range("B4").Value = "emlab"

# --- Sheet2: rename tab and populate with "data per year" data ------------
$ws2.Name = "data per year"

$ws2.Columns.Item(1).ColumnWidth = 23.09
$ws2.Columns.Item(2).ColumnWidth = 33.09
$ws2.Columns.Item(3).ColumnWidth = 23.09
$ws2.Columns.Item(4).ColumnWidth = 23.09
$ws2.Columns.Item(5).ColumnWidth = 23.09

$ws2.Range("A1").Value = "MarketClearingPoints"
$ws2.Range("B1").Value = "MarketClearingPoint 2021-08-04 13:27:00.769896"
$ws2.Range("C1").Value = "Market"
$ws2.Range("D1").Value = 1
$ws2.Range("E1").Value = "CO2Auction"

$ws2.Range("A2").Value = "MarketClearingPoints"
$ws2.Range("B2").Value = "MarketClearingPoint 2021-08-04 13:27:00.769896"
$ws2.Range("C2").Value = "Price"
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 0

$ws2.Range("A3").Value = "MarketClearingPoints"
$ws2.Range("B3").Value = "MarketClearingPoint 2021-08-04 13:27:00.769896"
$ws2.Range("C3").Value = "TotalCapacity"
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0

$ws2.Range("A4").Value = "PowerPlants"
$ws2.Range("B4").Value = "SunPV"
$ws2.Range("C4").Value = "MWNL"
$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = "12291,7"

$ws2.Range("A5").Value = "PowerPlants"
$ws2.Range("B5").Value = "WindOff"
$ws2.Range("C5").Value = "MWNL"
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = 3750

$ws2.Range("A6").Value = "PowerPlants"
$ws2.Range("B6").Value = "WindOn"
$ws2.Range("C6").Value = "MWNL"
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = 4590

$ws2.Range("A7").Value = "SystemClockTicks"
$ws2.Range("B7").Value = "SystemClockTicks"
$ws2.Range("C7").Value = "ticks"
$ws2.Range("D7").Value = 1
$ws2.Range("E7").Value = 1

$ws2.Range("A10").Value = "data per year in competes-emlab emlab"

# --- View state: keep Sheet1 active/selected, give Sheet2 its own selection
$ws2.Activate() | Out-Null
$ws2.Range("B14").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("N7").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 4
$win.ScrollColumn = 1
